# Generate Report for Handoff
# Updates the localization-status report to reflect that "b.md" has now
# been handed off (for both zh-cn and de-de locales) instead of already
# being handed back in sync with en-US.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Overview sheet: row for b.md (row 3) moves from "Handed back: in sync
# with en-US" to "Ready for handoff", with the new generation datetime.
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E3").Value = "Ready for handoff"
$overview.Range("F3").Value = "Ready for handoff"
$overview.Range("G3").Value = "2016-08-13 10:45:07"

# ---------------------------------------------------------------------
# zh-cn sheet: row for b.md (row 3) status + new handoff file/datetime
# and a new error detail explaining the handback file is stale.
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")
$zhcn.Range("C3").Value = "Ready for handoff"
$zhcn.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$zhcn.Range("H3").Value = "2016-08-13 10:44:55"
$zhcn.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/d52672629dbdaaa2c815f856a16cab1221386f38/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/498c9656d3af649e977ad5055e1b53436e635ea3/e2e/b.md."
$zhcn.Columns.Item(16).ColumnWidth = 39.17

# ---------------------------------------------------------------------
# de-de sheet: row for b.md (row 3) status + new handoff file/datetime
# and a new error detail explaining the handback file is stale.
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")
$dede.Range("C3").Value = "Ready for handoff"
$dede.Range("G3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$dede.Range("H3").Value = "2016-08-13 10:45:07"
$dede.Range("P3").Value = "The version of handback file is not the latest, current: https://github.com/OpenLocalizationTestOrg/oltest/blob/d52672629dbdaaa2c815f856a16cab1221386f38/e2e/a.md, latest: https://github.com/OpenLocalizationTestOrg/oltest/blob/498c9656d3af649e977ad5055e1b53436e635ea3/e2e/b.md."
$dede.Columns.Item(16).ColumnWidth = 39.17
